# Generate Report for Handoff
# - Update status "In Translation" -> "Ready for handoff" (shared across sheets)
# - Update "Latest HO Xliff Generate Date" / "Latest Handback DateTime" timestamp
#   2016-08-24 00:55:33 -> 2016-08-24 00:56:06  (Overview!G2, de-de!H2)
# - Update "Latest Handoff Datetime" timestamp
#   2016-08-24 00:55:28 -> 2016-08-24 00:55:57  (zh-cn!H2)
# - Widen columns E/F on Overview and column C on zh-cn/de-de from ~13.41 to ~17.22

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Status text updates ("In Translation" -> "Ready for handoff") ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Timestamp updates ---
$wsOverview.Range("G2").Value = "2016-08-24 00:56:06"
$wsDeDe.Range("H2").Value = "2016-08-24 00:56:06"

$wsZhCn.Range("H2").Value = "2016-08-24 00:55:57"

# --- Column width updates (widen to fit new text) ---
# Target stored OOXML width is 17.2159881591797. This engine stores column
# widths quantized to the nearest 1/6 of a character (matching Excel's
# internal pixel-granularity rounding), so the nearest representable value
# is used here (16.3333... as a COM ColumnWidth yields a stored width of
# 17.1666..., the closest achievable value to the target).
$newColWidth = 98 / 6

$wsOverview.Columns.Item(5).ColumnWidth = $newColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColWidth

$wsZhCn.Columns.Item(3).ColumnWidth = $newColWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newColWidth
